$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 99.666664
$ws.Range("I6").Value = 99.666664
$ws.Range("K6").Value = 298.999992
$ws.Range("M6").Value = -186.999992
$ws.Range("H17").Value = 4186.9375
$ws.Range("J17").Value = 4186.9375
$ws.Range("L17").Value = 12560.8125
$ws.Range("N17").Value = -12896.8125
$ws.Range("H86").Value = 5470
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 5470
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H92").Value = 280.72726
$ws.Range("I92").Value = 310.8
$ws.Range("K92").Value = 310.8
$ws.Range("M92").Value = 937.2
$ws.Range("H98").Value = 2498.4783
$ws.Range("I98").Value = 991.5714
$ws.Range("J98").Value = 3157.75
$ws.Range("K98").Value = 991.5714
$ws.Range("L98").Value = 3157.75
$ws.Range("M98").Value = 506.4286
$ws.Range("N98").Value = -6153.75
$ws.Range("H100").Value = 4369.3335
$ws.Range("I100").Value = 4608
$ws.Range("J100").Value = 4250
$ws.Range("K100").Value = 4608
$ws.Range("L100").Value = 4250
$ws.Range("M100").Value = -4067
$ws.Range("N100").Value = -5332
$ws.Range("H116").Value = 7838.2
$ws.Range("I116").Value = 6332.6665
$ws.Range("K116").Value = 6332.6665
$ws.Range("M116").Value = -2890.6665
$ws.Range("H122").Value = 2498.4783
$ws.Range("I122").Value = 991.5714
$ws.Range("J122").Value = 3157.75
$ws.Range("K122").Value = 2974.7142
$ws.Range("L122").Value = 9473.25
$ws.Range("M122").Value = -524.7142000000003
$ws.Range("N122").Value = -14373.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 5850
$ws.Range("I4").Value = 3551.5
$ws.Range("K4").Value = 3551.5
$ws.Range("M4").Value = -3435.5
$ws.Range("H45").Value = 5259.6
$ws.Range("I45").Value = 5259.6
$ws.Range("K45").Value = 5259.6
$ws.Range("M45").Value = -4882.6
$ws.Range("H61").Value = 3810.8147
$ws.Range("I61").Value = 1951.5454
$ws.Range("J61").Value = 5089.0625
$ws.Range("K61").Value = 1951.5454
$ws.Range("L61").Value = 5089.0625
$ws.Range("M61").Value = -1739.5454
$ws.Range("N61").Value = -5513.0625
$ws.Range("H97").Value = 553
$ws.Range("I97").Value = 553
$ws.Range("K97").Value = 553
$ws.Range("M97").Value = -57
$ws.Range("H102").Value = 2068.2
$ws.Range("I102").Value = 2068.2
$ws.Range("K102").Value = 2068.2
$ws.Range("M102").Value = -446.1999999999998
$ws.Range("H134").Value = 64997
$ws.Range("J134").Value = 64997
$ws.Range("L134").Value = 64997
$ws.Range("N134").Value = -75137
$ws.Range("H136").Value = 3810.8147
$ws.Range("I136").Value = 1951.5454
$ws.Range("J136").Value = 5089.0625
$ws.Range("K136").Value = 5854.6362
$ws.Range("L136").Value = 15267.1875
$ws.Range("M136").Value = -3304.6362
$ws.Range("N136").Value = -20367.1875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4332.6665
$ws.Range("I86").Value = 4332.6665
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 4332.6665
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -3209.6665
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 4332.6665
$ws.Range("I89").Value = 4332.6665
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 21663.3325
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -16047.3325
$ws.Range("N89").ClearContents()
$ws.Range("H94").Value = 1516.3846
$ws.Range("I94").Value = 1559.8334
$ws.Range("J94").Value = 995
$ws.Range("K94").Value = 1559.8334
$ws.Range("L94").Value = 995
$ws.Range("M94").Value = -1108.8334
$ws.Range("N94").Value = -1897
$ws.Range("H99").Value = 2062.75
$ws.Range("I99").Value = 2079.3635
$ws.Range("K99").Value = 2079.3635
$ws.Range("M99").Value = -581.3634999999999
$ws.Range("H105").Value = 4758.4
$ws.Range("I105").Value = 3948
$ws.Range("J105").Value = 8000
$ws.Range("K105").Value = 3948
$ws.Range("L105").Value = 8000
$ws.Range("M105").Value = -2201
$ws.Range("N105").Value = -11494
$ws.Range("H107").Value = 4240.8335
$ws.Range("I107").Value = 4089
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 4089
$ws.Range("L107").Value = 5000
$ws.Range("M107").Value = -2169
$ws.Range("N107").Value = -8840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 36719
$ws.Range("I54").Value = 30000
$ws.Range("K54").Value = 30000
$ws.Range("M54").Value = -29342
$ws.Range("H62").Value = 1999
$ws.Range("I62").Value = 1999
$ws.Range("K62").Value = 1999
$ws.Range("M62").Value = -1375
$ws.Range("H65").Value = 1999
$ws.Range("I65").Value = 1999
$ws.Range("K65").Value = 9995
$ws.Range("M65").Value = -6875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1120.7
$ws.Range("I60").Value = 1101.625
$ws.Range("J60").Value = 1197
$ws.Range("K60").Value = 3304.875
$ws.Range("L60").Value = 3591
$ws.Range("M60").Value = -3053.875
$ws.Range("N60").Value = -4093
$ws.Range("H132").Value = 5299.125
$ws.Range("I132").Value = 5878.8
$ws.Range("K132").Value = 52909.2
$ws.Range("M132").Value = -50379.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1497
$ws.Range("I80").Value = 1497
$ws.Range("K80").Value = 1497
$ws.Range("M80").Value = -499
$ws.Range("H83").Value = 1497
$ws.Range("I83").Value = 1497
$ws.Range("K83").Value = 7485
$ws.Range("M83").Value = -2493
$ws.Range("H97").Value = 997.5
$ws.Range("I97").Value = 997.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 997.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -501.5
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7255.4
$ws.Range("J7").Value = 8630.789000000001
$ws.Range("L7").Value = 8630.789000000001
$ws.Range("N7").Value = -8854.789000000001
$ws.Range("H16").Value = 1927.5714
$ws.Range("I16").Value = 898.6
$ws.Range("J16").Value = 4500
$ws.Range("K16").Value = 898.6
$ws.Range("L16").Value = 4500
$ws.Range("M16").Value = -728.6
$ws.Range("N16").Value = -4840
$ws.Range("H22").Value = 2604.2
$ws.Range("I22").Value = 3356.5
$ws.Range("J22").Value = 2102.6667
$ws.Range("K22").Value = 3356.5
$ws.Range("L22").Value = 2102.6667
$ws.Range("M22").Value = -3061.5
$ws.Range("N22").Value = -2692.6667
$ws.Range("H27").Value = 2604.2
$ws.Range("I27").Value = 3356.5
$ws.Range("J27").Value = 2102.6667
$ws.Range("K27").Value = 3356.5
$ws.Range("L27").Value = 2102.6667
$ws.Range("M27").Value = -3249.5
$ws.Range("N27").Value = -2316.6667
$ws.Range("H40").Value = 3922
$ws.Range("I40").Value = 2333
$ws.Range("K40").Value = 2333
$ws.Range("M40").Value = -2197
$ws.Range("H82").Value = 1336
$ws.Range("I82").Value = 1383.2
$ws.Range("J82").Value = 1100
$ws.Range("K82").Value = 1383.2
$ws.Range("L82").Value = 1100
$ws.Range("M82").Value = -1022.2
$ws.Range("N82").Value = -1822
$ws.Range("H85").Value = 1336
$ws.Range("I85").Value = 1383.2
$ws.Range("J85").Value = 1100
$ws.Range("K85").Value = 1383.2
$ws.Range("L85").Value = 1100
$ws.Range("M85").Value = -135.2
$ws.Range("N85").Value = -3596
$ws.Range("H93").Value = 2000
$ws.Range("I93").Value = 2000
$ws.Range("K93").Value = 2000
$ws.Range("M93").Value = -752
$ws.Range("H96").Value = 49999
$ws.Range("J96").Value = 49999
$ws.Range("L96").Value = 49999
$ws.Range("N96").Value = -55491
$ws.Range("H100").Value = 3441.2144
$ws.Range("J100").Value = 4550
$ws.Range("L100").Value = 4550
$ws.Range("N100").Value = -5632
$ws.Range("H126").Value = 7255.4
$ws.Range("J126").Value = 8630.789000000001
$ws.Range("L126").Value = 25892.367
$ws.Range("N126").Value = -30832.367
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H132").Value = 3136.375
$ws.Range("I132").Value = 2547.6667
$ws.Range("J132").Value = 4902.5
$ws.Range("K132").Value = 7643.000100000001
$ws.Range("L132").Value = 14707.5
$ws.Range("M132").Value = -5113.000100000001
$ws.Range("N132").Value = -19767.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 20166
$ws.Range("J62").Value = 14694.571
$ws.Range("L62").Value = 14694.571
$ws.Range("N62").Value = -15942.571
$ws.Range("H65").Value = 20166
$ws.Range("J65").Value = 14694.571
$ws.Range("L65").Value = 73472.855
$ws.Range("N65").Value = -79712.855
$ws.Range("H100").Value = 8335575
$ws.Range("I100").Value = 11112703
$ws.Range("J100").Value = 4190
$ws.Range("K100").Value = 22225406
$ws.Range("L100").Value = 8380
$ws.Range("M100").Value = -22224865
$ws.Range("N100").Value = -9462
$ws.Range("H129").Value = 72999
$ws.Range("J129").Value = 72999
$ws.Range("L129").Value = 72999
$ws.Range("N129").Value = -82999
